$wb = $excel.ActiveWorkbook

# --- ProductDetail: insert a new "ProdDetailDesc / string" row (row 4) ---
$wsProductDetail = $wb.Worksheets.Item("ProductDetail")
$wsProductDetail.Rows.Item(4).Insert()
$wsProductDetail.Range("B4").Value = "ProdDetailDesc"
$wsProductDetail.Range("C4").Value = "string"
$wsProductDetail.Columns.Item(2).ColumnWidth = 12.65

# --- OrderDetail: insert a new "FK / OrderId / int" row (row 3) ---
$wsOrderDetail = $wb.Worksheets.Item("OrderDetail")
$wsOrderDetail.Rows.Item(3).Insert()
$wsOrderDetail.Range("A3").Value = "FK"
$wsOrderDetail.Range("B3").Value = "OrderId"
$wsOrderDetail.Range("C3").Value = "int"

# --- Remove the stray, unused "Sheet6" worksheet ---
$wsStray = $wb.Worksheets.Item("Sheet6")
$wsStray.Delete()

# --- Update selections on the various sheets ---
$wb.Worksheets.Item("Product").Range("D6").Select()
$wb.Worksheets.Item("Order").Range("B2:C2").Select()
$wsOrderDetail.Range("F17").Select()
$wb.Worksheets.Item("User").Range("L23").Select()

# Select ProductDetail last so it becomes the active sheet/tab (activeTab=1)
$wsProductDetail.Range("C5").Select()
